$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formula edits: B column changes from SIN(x) to SIN(2*x) ---
# B2 is a standalone formula (not part of either shared-formula group).
$ws.Range("B2").Formula = "=SIN(2*A2)"

# B3:B66 is the first shared-formula group (anchor B3).
$ws.Range("B3:B66").Formula = "=SIN(2*A3)"

# B67:B102 is the second shared-formula group (anchor B67).
$ws.Range("B67:B102").Formula = "=SIN(2*A67)"

# --- View/selection state ---
$ws.Range("F11").Select()
